$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the dSF column (F) values for the rows that changed in the repull
$ws.Range("F2").Value  = 3
$ws.Range("F3").Value  = -3
$ws.Range("F6").Value  = 3
$ws.Range("F7").Value  = -2
$ws.Range("F8").Value  = 7
$ws.Range("F13").Value = 5
$ws.Range("F15").Value = 2
$ws.Range("F21").Value = -3
$ws.Range("F24").Value = 5
